# Supplier Excel template: add a "SupplierType" column and a second
# example supplier row, matching the updated import/reporting feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "SupplierType" header + the value for the existing sample row
$ws.Range("L1").Value = "SupplierType"
$ws.Range("L2").Value = "Individual"

# New second example supplier (row 3) mirroring the sample row's data,
# except for the supplier name
$ws.Range("A3").Value = "Dami and Sons"
$ws.Range("B3").Value = "johndoe@gmail.com"
$ws.Range("C3").Value = 542542299
$ws.Range("E3").Value = "Accra"
$ws.Range("F3").Value = "GD-898-0909"
$ws.Range("G3").Value = "Minoxidil"
$ws.Range("H3").Value = "1 month"
$ws.Range("L3").Value = "Company"

# Mirror the hyperlink + hyperlink styling from the existing B2 sample cell
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:johndoe@gmail.com")
$ws.Range("B3").Style = "Hyperlink"

# Restore the workbook's saved view/selection state
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L3").Select() | Out-Null
